$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row above row 15 (current "[MKS SBASE v1.3]" entry),
# shifting all subsequent rows down by one.
$ws.Rows.Item(15).Insert()

# Fill in the new controller/driver pair.
$ws.Range("A15").Value = "[Mach 3 USB BOB - BSMCE04U-PP](https://drufelcnc.com/?c=controllers&p=BSMCE04U)<sup>1</sup>"
$ws.Range("B15").Value = "[STM32F1xx](https://github.com/grblHAL/STM32F1xx)"

# Update the active selection to match the saved view state.
$ws.Range("A32").Select()
